$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the sheet view (scroll position + active selection) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1

# --- Append the new data row (row 65), mirroring row 64's values with the date advanced one day ---
$ws.Range("A65").Value = 43569
$ws.Range("B65").Value = 3
$ws.Range("C65").Value = 57
$ws.Range("D65").Value = 73
$ws.Range("E65").Value = 80
$ws.Range("F65").Value = 56
$ws.Range("G65").Value = 14
$ws.Range("H65").Value = 5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 77
$ws.Range("K65").Value = 35
$ws.Range("L65").Value = 35
$ws.Range("M65").Value = 25
$ws.Range("N65").Value = 12
$ws.Range("O65").Value = 4

# --- Finally move/select the active cell to match the saved view state ---
$ws.Range("A62").Select()
